# This script applies updated market-price data (columns H:N) produced by
# the scheduled pricing runner to the relevant rows on each Leve sheet.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 20013
$ws.Range("I26").Value = 20013
$ws.Range("K26").Value = 20013
$ws.Range("M26").Value = -19669
# Row 103
$ws.Range("H103").Value = 799.8333
$ws.Range("I103").Value = 833
$ws.Range("K103").Value = 2499
$ws.Range("M103").Value = -1913
# Row 112
$ws.Range("H112").Value = 2808.889
$ws.Range("J112").Value = 2808.889
$ws.Range("L112").Value = 8426.667000000001
$ws.Range("N112").Value = -10642.667
# Row 129
$ws.Range("H129").Value = 1009.6964
$ws.Range("I129").Value = 596.6667
$ws.Range("J129").Value = 1033.0754
$ws.Range("K129").Value = 1790.0001
$ws.Range("L129").Value = 3099.2262
$ws.Range("M129").Value = 3209.9999
$ws.Range("N129").Value = -13099.2262
# Row 137
$ws.Range("H137").Value = 1650.6538
$ws.Range("I137").Value = 1374.3334
$ws.Range("J137").Value = 2027.4546
$ws.Range("K137").Value = 4123.0002
$ws.Range("L137").Value = 6082.3638
$ws.Range("M137").Value = -1573.0002
$ws.Range("N137").Value = -11182.3638

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1181.6666
$ws.Range("I61").Value = 1181.6666
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1181.6666
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -969.6666
$ws.Range("N61").Value = $null
# Row 102
$ws.Range("H102").Value = 2540.125
$ws.Range("I102").Value = 2282
$ws.Range("K102").Value = 2282
$ws.Range("M102").Value = -660
# Row 136
$ws.Range("H136").Value = 1181.6666
$ws.Range("I136").Value = 1181.6666
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3544.9998
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -994.9998000000001
$ws.Range("N136").Value = $null

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 2469327.8
$ws.Range("I80").Value = 12345745
$ws.Range("J80").Value = 223.33333
$ws.Range("K80").Value = 12345745
$ws.Range("L80").Value = 223.33333
$ws.Range("M80").Value = -12344747
$ws.Range("N80").Value = -2219.33333
# Row 83
$ws.Range("H83").Value = 2469327.8
$ws.Range("I83").Value = 12345745
$ws.Range("J83").Value = 223.33333
$ws.Range("K83").Value = 61728725
$ws.Range("L83").Value = 1116.66665
$ws.Range("M83").Value = -61723733
$ws.Range("N83").Value = -11100.66665
# Row 129
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2234.64
$ws.Range("I31").Value = 1668.1
$ws.Range("K31").Value = 1668.1
$ws.Range("M31").Value = -1373.1
# Row 32
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -684
# Row 34
$ws.Range("H34").Value = 2234.64
$ws.Range("I34").Value = 1668.1
$ws.Range("K34").Value = 1668.1
$ws.Range("M34").Value = -1466.1
# Row 58
$ws.Range("H58").Value = 1244.4117
$ws.Range("I58").Value = 1400.2727
$ws.Range("J58").Value = 958.6667
$ws.Range("K58").Value = 1400.2727
$ws.Range("L58").Value = 958.6667
$ws.Range("M58").Value = -1197.2727
$ws.Range("N58").Value = -1364.6667
# Row 99
$ws.Range("H99").Value = 3206.353
$ws.Range("I99").Value = 3391.0833
$ws.Range("J99").Value = 2763
$ws.Range("K99").Value = 3391.0833
$ws.Range("L99").Value = 2763
$ws.Range("M99").Value = -1893.0833
$ws.Range("N99").Value = -5759
# Row 105
$ws.Range("H105").Value = 1105.5
$ws.Range("I105").Value = 1105.5
$ws.Range("K105").Value = 1105.5
$ws.Range("M105").Value = 641.5
# Row 126
$ws.Range("H126").Value = 3206.353
$ws.Range("I126").Value = 3391.0833
$ws.Range("J126").Value = 2763
$ws.Range("K126").Value = 10173.2499
$ws.Range("L126").Value = 8289
$ws.Range("M126").Value = -7703.249899999999
$ws.Range("N126").Value = -13229
# Row 136
$ws.Range("H136").Value = 1244.4117
$ws.Range("I136").Value = 1400.2727
$ws.Range("J136").Value = 958.6667
$ws.Range("K136").Value = 4200.8181
$ws.Range("L136").Value = 2876.0001
$ws.Range("M136").Value = -1650.8181
$ws.Range("N136").Value = -7976.0001

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1517.5714
$ws.Range("I4").Value = 61.5
$ws.Range("J4").Value = 2100
$ws.Range("K4").Value = 184.5
$ws.Range("L4").Value = 6300
$ws.Range("M4").Value = -72.5
$ws.Range("N4").Value = -6524
# Row 5
$ws.Range("H5").Value = 1463.6666
$ws.Range("I5").Value = 1642.5294
$ws.Range("K5").Value = 4927.5882
$ws.Range("M5").Value = -4815.5882
# Row 41
$ws.Range("J41").Value = 1020
$ws.Range("L41").Value = 3060
$ws.Range("N41").Value = -3736
# Row 129
$ws.Range("H129").Value = 2779263.5
$ws.Range("I129").Value = 907.25
$ws.Range("J129").Value = 3573079.8
$ws.Range("K129").Value = 2721.75
$ws.Range("L129").Value = 10719239.4
$ws.Range("M129").Value = 2278.25
$ws.Range("N129").Value = -10729239.4
# Row 135
$ws.Range("H135").Value = 1463.6666
$ws.Range("I135").Value = 1642.5294
$ws.Range("K135").Value = 14782.7646
$ws.Range("M135").Value = -12247.7646

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 49996.668
$ws.Range("I35").Value = 49995
$ws.Range("J35").Value = 50000
$ws.Range("K35").Value = 49995
$ws.Range("L35").Value = 50000
$ws.Range("M35").Value = -49697
$ws.Range("N35").Value = -50596
# Row 97
$ws.Range("H97").Value = 39725.383
$ws.Range("I97").Value = 46720
$ws.Range("J97").Value = 1255
$ws.Range("K97").Value = 46720
$ws.Range("L97").Value = 1255
$ws.Range("M97").Value = -46224
$ws.Range("N97").Value = -2247
# Row 109
$ws.Range("H109").Value = 20285
$ws.Range("J109").Value = 20285
$ws.Range("L109").Value = 20285
$ws.Range("N109").Value = -22365
# Row 122
$ws.Range("H122").Value = 1466.909
$ws.Range("I122").Value = 1202
$ws.Range("J122").Value = 2367.6
$ws.Range("K122").Value = 3606
$ws.Range("L122").Value = 7102.799999999999
$ws.Range("M122").Value = -1156
$ws.Range("N122").Value = -12002.8

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1025
$ws.Range("I22").Value = 485.7143
$ws.Range("K22").Value = 485.7143
$ws.Range("M22").Value = -190.7143
# Row 27
$ws.Range("H27").Value = 1025
$ws.Range("I27").Value = 485.7143
$ws.Range("K27").Value = 485.7143
$ws.Range("M27").Value = -378.7143

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 25418.234
$ws.Range("J123").Value = 25418.234
$ws.Range("L123").Value = 25418.234
$ws.Range("N123").Value = -35218.234
# Row 126
$ws.Range("H126").Value = 9892.556
$ws.Range("I126").Value = 9892.556
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 29677.668
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27207.668
$ws.Range("N126").Value = $null
